$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New values for the 100 math-fact cells, in row-major order (20 rows x 5 cols)
$newValues = @(
"18+7=", "16+65=", "76-44=", "59-21=", "53+13=", "30-14=", "48+46=", "6+18=", "3+55=", "17+19=", "42+53=", "57+33=", "70+11=", "72+14=", "88-9=", "65-33=", "63+32=", "0+88=", "28+29=", "71-58=", "75-41=", "39+8=", "69-17=", "99-82=", "35+33=", "30+1=", "82-64=", "80-46=", "53-37=", "46+29=", "41+9=", "52+0=", "23-18=", "8+70=", "14+28=", "86-4=", "71-15=", "46-3=", "87-57=", "3+83=", "44+43=", "19+20=", "87-18=", "80-19=", "60-40=", "85-28=", "22+21=", "41+23=", "19+27=", "65-29=", "91-56=", "50+21=", "50+17=", "18+24=", "31+53=", "17+61=", "69+3=", "20+3=", "34+58=", "13+76=", "1+73=", "53-7=", "63+26=", "74+2=", "89+5=", "25-14=", "62+7=", "72-1=", "54-30=", "62-32=", "51+45=", "8+85=", "3+90=", "15+7=", "30-17=", "72-65=", "17-11=", "43+26=", "14+1=", "19-16=", "10+61=", "88-82=", "86-0=", "62+27=", "15+50=", "59-35=", "88-44=", "50-41=", "25-0=", "68+16=", "53-32=", "44+21=", "88-72=", "54+18=", "20+19=", "12+21=", "80+7=", "61+12=", "93-18=", "59-26="
)

$rows = 20
$cols = 5
for ($i = 0; $i -lt ($rows * $cols); $i++) {
    $row = [math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 2
    $r.Text = $newValues[$i]
}
